$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("9:9").Insert() | Out-Null
$ws.Range("B10:E10").Copy() | Out-Null
$ws.Range("B9:E9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("B9").Value = "boolean"
$ws.Range("C9").Value = """type"": ""boolean"""
$ws.Range("C9").Select() | Out-Null
